$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.236445426940918
$ws.Range("B1").Value = 6.718282699584961
$ws.Range("C1").Value = 6.538310050964355
$ws.Range("D1").Value = 2.298135042190552
$ws.Range("E1").Value = 1.482321619987488
